# Apply updated cryptos list values to worksheet (scraped Thu Jul  4 20:13:21 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '58.190.77'
$ws.Cells.Item(2, 5).Value = '  -2.63%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.129.81'
$ws.Cells.Item(3, 5).Value = '  -4.45%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '523.57'
$ws.Cells.Item(5, 5).Value = '  -5.20%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '133.40'
$ws.Cells.Item(6, 5).Value = '  -4.93%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.03%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '3.128.03'
$ws.Cells.Item(8, 5).Value = '  -4.53%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.441'
$ws.Cells.Item(9, 5).Value = '  -4.79%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '7.20'
$ws.Cells.Item(10, 5).Value = '  -7.47%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -8.51%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.377'
$ws.Cells.Item(12, 5).Value = '  -6.61%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '3.666.90'
$ws.Cells.Item(13, 5).Value = '  -4.42%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  -1.11%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '25.45'
$ws.Cells.Item(15, 5).Value = '  -4.14%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '3.129.91'
$ws.Cells.Item(16, 5).Value = '  -4.28%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '58.166.76'
$ws.Cells.Item(17, 5).Value = '  -2.90%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  -6.60%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '5.72'
$ws.Cells.Item(19, 5).Value = '  -5.58%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '12.89'
$ws.Cells.Item(20, 5).Value = '  -5.65%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '7.84'
$ws.Cells.Item(21, 5).Value = '  -7.73%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '342.88'
$ws.Cells.Item(22, 5).Value = '  -7.79%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  -0.05%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '67.87'
$ws.Cells.Item(24, 5).Value = '  -7.43%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.506'
$ws.Cells.Item(25, 5).Value = '  -4.43%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '3.265.91'
$ws.Cells.Item(26, 5).Value = '  -4.24%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.169'
$ws.Cells.Item(27, 5).Value = '  +0.09%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '1.00'
$ws.Cells.Item(28, 5).Value = '  -0.04%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'PEPE'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '0.0₃0946'
$ws.Cells.Item(29, 5).Value = '  -6.19%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '6.78'
$ws.Cells.Item(30, 5).Value = '  -3.87%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.998'
$ws.Cells.Item(31, 5).Value = '  -0.14%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -7.49%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  -7.63%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.22'
$ws.Cells.Item(34, 5).Value = '  -0.74%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '21.26'
$ws.Cells.Item(35, 5).Value = '  -5.15%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'Monero'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '158.16'
$ws.Cells.Item(36, 5).Value = '  -4.75%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'NEARProtocol'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '4.78'
$ws.Cells.Item(37, 5).Value = '  -4.99%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -5.77%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.36'
$ws.Cells.Item(39, 5).Value = '  -10.21%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'Hedera'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.0683'
$ws.Cells.Item(40, 5).Value = '  -5.53%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'RenzoRestakedETH'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '3.160.48'
$ws.Cells.Item(41, 5).Value = '  -4.42%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '40.44'
$ws.Cells.Item(42, 5).Value = '  -2.96%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '23.82'
$ws.Cells.Item(43, 5).Value = '  -8.11%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.689'
$ws.Cells.Item(44, 5).Value = '  -7.28%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -1.63%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '3.88'
$ws.Cells.Item(46, 5).Value = '  -5.22%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -0.08%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'Maker'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.271.60'
$ws.Cells.Item(48, 5).Value = '  -2.20%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Stacks'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.43'
$ws.Cells.Item(49, 5).Value = '  -7.94%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '6.16'
$ws.Cells.Item(50, 5).Value = '  -2.38%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '20.67'
$ws.Cells.Item(51, 5).Value = '  -1.93%  '
